# Apply the "Deploying to gh-pages ... LinuxForHealth/alvearie-fhir-ig" edit
# to the StructureDefinition-cdm-identifier workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Metadata" sheet - top level StructureDefinition metadata
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL: ibm.com -> linuxforhealth.org
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/cdm-identifier"

# Version: 7.0.0 -> 8.0.0
$meta.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---------------------------------------------------------------------------
# "Elements" sheet - element definition table
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Identifier" element. Its Constraint(s) cell (column AI)
# is cleared.
$elements.Range("AI2").Value = ""

# Row 5 is the "encryptedState" extension slice (Identifier.extension).
# Its Type(s) cell (column J) references the encrypted-state extension URL,
# which moves from ibm.com to linuxforhealth.org.
$elements.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/encrypted-state}" + [char]10
